# Apply "Update countries & provincias Spain" edits to the Pais sheet.
# The shared-string table order shifted (countries inserted/reordered) and
# several case-count statistics were refreshed; this sets every affected
# cell to its final value directly (row/column position is unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 17:46"

$ws.Range("B8").Value = 18756
$ws.Range("C8").Value = 3436
$ws.Range("E8").Value = 18523
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 53

$ws.Range("B9").Value = 16058
$ws.Range("C9").Value = 2269
$ws.Range("E9").Value = 15714

$ws.Range("F14").Value = 210

$ws.Range("B17").Value = 1906
$ws.Range("C17").Value = 116
$ws.Range("E17").Value = 1898

$ws.Range("A30").Value = "Pakistan"
$ws.Range("B30").Value = 500
$ws.Range("C30").Value = 46
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = 484
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 3

$ws.Range("A31").Value = "Grecia"
$ws.Range("B31").Value = 495
$ws.Range("C31").Value = 31
$ws.Range("D31").Value = 19
$ws.Range("E31").Value = 467
$ws.Range("F31").Value = 20
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 9

$ws.Range("A32").Value = "Luxemburgo"
$ws.Range("B32").Value = 484
$ws.Range("C32").Value = 149
$ws.Range("D32").Value = 6
$ws.Range("E32").Value = 473
$ws.Range("F32").Value = 1
$ws.Range("H32").Value = 5

$ws.Range("A52").Value = "India"
$ws.Range("B52").Value = 234
$ws.Range("C52").Value = 40
$ws.Range("D52").Value = 23
$ws.Range("E52").Value = 206
$ws.Range("F52").Value = 0
$ws.Range("H52").Value = 5

$ws.Range("A53").Value = "Filipinas"
$ws.Range("B53").Value = 230
$ws.Range("C53").Value = 13
$ws.Range("D53").Value = 8
$ws.Range("E53").Value = 204
$ws.Range("F53").Value = 1
$ws.Range("H53").Value = 18

$ws.Range("A78").Value = "Republica de Chipre"
$ws.Range("C78").Value = 8
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 75
$ws.Range("F78").Value = 1

$ws.Range("A79").Value = "Principado de Andorra"
$ws.Range("B79").Value = 75
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 1
$ws.Range("E79").Value = 74
$ws.Range("F79").Value = 2
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 0

$ws.Range("A80").Value = "Marruecos"
$ws.Range("B80").Value = 74
$ws.Range("C80").Value = 11
$ws.Range("D80").Value = 2
$ws.Range("E80").Value = 69
$ws.Range("F80").Value = 1
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 3

$ws.Range("A81").Value = "Sri Lanka"
$ws.Range("B81").Value = 73
$ws.Range("C81").Value = 13
$ws.Range("D81").Value = 3
$ws.Range("H81").Value = 0

$ws.Range("A82").Value = "Republica Dominicana"
$ws.Range("B82").Value = 72
$ws.Range("C82").Value = 38
$ws.Range("D82").Value = 0
$ws.Range("E82").Value = 70
$ws.Range("F82").Value = 0
$ws.Range("H82").Value = 2

$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("C83").Value = 20
$ws.Range("D83").Value = 1
$ws.Range("E83").Value = 69
$ws.Range("F83").Value = 1
$ws.Range("H83").Value = 0

$ws.Range("A84").Value = "Albania"
$ws.Range("B84").Value = 70
$ws.Range("C84").Value = 6
$ws.Range("D84").Value = 0
$ws.Range("F84").Value = 2
$ws.Range("H84").Value = 2

$ws.Range("A85").Value = "Jordania"
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 1
$ws.Range("E85").Value = 68
$ws.Range("F85").Value = 0

$ws.Range("A86").Value = "Bosnia y Herzegovina"
$ws.Range("C86").Value = 5
$ws.Range("D86").Value = 2
$ws.Range("E86").Value = 67
$ws.Range("F86").Value = 1

$ws.Range("A87").Value = "Bielorrusia"
$ws.Range("B87").Value = 69
$ws.Range("C87").Value = 18
$ws.Range("D87").Value = 15
$ws.Range("E87").Value = 54
$ws.Range("F87").Value = 0

$ws.Range("A91").Value = "Kazajistan"
$ws.Range("B91").Value = 52
$ws.Range("C91").Value = 8
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 52

$ws.Range("A92").Value = "Camboya"
$ws.Range("B92").Value = 51
$ws.Range("C92").Value = 14
$ws.Range("D92").Value = 1
$ws.Range("E92").Value = 50

$ws.Range("A121").Value = "Maldivas"

$ws.Range("A123").Value = "Montenegro"

$ws.Range("A166").Value = "Isla de Man"
$ws.Range("C166").Value = 1

$ws.Range("A167").Value = "Haiti"
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 2
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = "Sudan"
$ws.Range("B168").Value = 2
$ws.Range("H168").Value = 1

$ws.Range("A171").Value = "San Martin (Parte Holandesa)"

$ws.Range("A172").Value = "Antigua y Barbuda"

$ws.Range("A173").Value = "El Salvador"
$ws.Range("C173").Value = 0

$ws.Range("A174").Value = "Nicaragua"

$ws.Range("A175").Value = "Republica del Chad"

$ws.Range("A176").Value = "Papua Nueva Guinea"
$ws.Range("C176").Value = 1

$ws.Range("A177").Value = "Montserrat"
$ws.Range("C177").Value = 0

$ws.Range("A178").Value = "Gambia"

$ws.Range("A179").Value = "Niger"

$ws.Range("A180").Value = "Santa Sede"

$ws.Range("A181").Value = "Fiyi"

$ws.Range("A182").Value = "Cabo Verde"
$ws.Range("C182").Value = 1

$ws.Range("A183").Value = "San Vicente y las Granadinas"

$ws.Range("A184").Value = "Somalia"

$ws.Range("A185").Value = "Republica de Yibuti"
